# Updated cryptos list on Wed Feb 22 21:35:59 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must stay text (Price column has numeric-looking
# strings like "308.32" that Excel's COM layer would otherwise auto-convert
# to a real number). Force text via NumberFormat, assign, then clear the
# temporary formatting so the cell's style index is untouched (matches the
# source workbook, where these cells carry no explicit style).
function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

# Row 2 - Bitcoin
Set-TextValue "D2" "23.821.57"
$ws.Range("E2").Value = "  -2.33%  "

# Row 3 - Ethereum
$ws.Range("E3").Value = "  -2.01%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.65%  "

# Row 5 - BNB
Set-TextValue "D5" "308.01"
$ws.Range("E5").Value = "  -1.30%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  -0.59%  "

# Row 7 - XRP
Set-TextValue "D7" "0.3938"
$ws.Range("E7").Value = "  +0.55%  "

# Row 8 - Cardano
Set-TextValue "D8" "0.3836"
$ws.Range("E8").Value = "  -1.60%  "

# Row 9 - BinanceUSD
$ws.Range("E9").Value = "  -0.44%  "

# Row 10 - was Polygon, now OKB
$ws.Range("B10").Value = "OKB"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D10" "49.29"
$ws.Range("E10").Value = "  -2.31%  "

# Row 11 - was OKB, now Polygon
$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D11" "1.351"
$ws.Range("E11").Value = "  -2.09%  "

# Row 12 - Dogecoin
Set-TextValue "D12" "0.08454"
$ws.Range("E12").Value = "  -0.98%  "

# Row 13 - Solana
Set-TextValue "D13" "23.67"
$ws.Range("E13").Value = "  -5.12%  "

# Row 14 - Polkadot
Set-TextValue "D14" "7.048"
$ws.Range("E14").Value = "  -2.37%  "

# Row 15 - Chainlink
Set-TextValue "D15" "7.566"
$ws.Range("E15").Value = "  -0.31%  "

# Row 16 - ShibaInu
Set-TextValue "D16" "0.00001282"
$ws.Range("E16").Value = "  -1.52%  "

# Row 17 - WrappedEther
Set-TextValue "D17" "1.614.87"
$ws.Range("E17").Value = "  -3.51%  "

# Row 18 - Litecoin
Set-TextValue "D18" "93.72"
$ws.Range("E18").Value = "  +0.70%  "

# Row 19 - TRON
Set-TextValue "D19" "0.06927"
$ws.Range("E19").Value = "  -0.63%  "

# Row 20 - Avalanche
Set-TextValue "D20" "19.93"
$ws.Range("E20").Value = "  -6.28%  "

# Row 21 - Uniswap
Set-TextValue "D21" "6.813"
$ws.Range("E21").Value = "  -2.38%  "

# Row 22 - Dai
Set-TextValue "D22" "1.000"
$ws.Range("E22").Value = "  -0.64%  "

# Row 23 - Cosmos
Set-TextValue "D23" "13.43"

# Row 24 - WrappedBTC
Set-TextValue "D24" "23.832.00"
$ws.Range("E24").Value = "  -2.28%  "

# Row 25 - Toncoin
Set-TextValue "D25" "2.477"
$ws.Range("E25").Value = "  +6.03%  "

# Row 26 - LidoDAOToken
Set-TextValue "D26" "2.830"
$ws.Range("E26").Value = "  +2.57%  "

# Row 27 - EthereumClassic
Set-TextValue "D27" "22.22"
$ws.Range("E27").Value = "  -2.19%  "

# Row 28 - Monero
Set-TextValue "D28" "156.94"
$ws.Range("E28").Value = "  -0.84%  "

# Row 29 - BitcoinCash
Set-TextValue "D29" "140.15"
$ws.Range("E29").Value = "  -2.96%  "

# Row 30 - HuobiToken
Set-TextValue "D30" "5.293"
$ws.Range("E30").Value = "  -8.48%  "

# Row 31 - Filecoin
Set-TextValue "D31" "7.808"
$ws.Range("E31").Value = "  -5.54%  "

# Row 32 - WEMIXTOKEN
Set-TextValue "D32" "2.488"
$ws.Range("E32").Value = "  -1.17%  "

# Row 33 - WrappedliquidstakedEther2.0
Set-TextValue "D33" "1.788.83"
$ws.Range("E33").Value = "  -3.11%  "

# Row 34 - Hedera
Set-TextValue "D34" "0.08112"
$ws.Range("E34").Value = "  -0.73%  "

# Row 35 - ImmutableX
Set-TextValue "D35" "0.9835"
$ws.Range("E35").Value = "  -1.16%  "

# Row 36 - InternetComputer(DFINITY)
Set-TextValue "D36" "6.616"
$ws.Range("E36").Value = "  -3.57%  "

# Row 37 - VeChain
Set-TextValue "D37" "0.02881"
$ws.Range("E37").Value = "  -4.49%  "

# Row 38 - Algorand
Set-TextValue "D38" "0.2669"
$ws.Range("E38").Value = "  -3.46%  "

# Row 39 - Stellar
Set-TextValue "D39" "0.09146"
$ws.Range("E39").Value = "  -4.12%  "

# Row 40 - FraxShare
Set-TextValue "D40" "10.37"
$ws.Range("E40").Value = "  +1.77%  "

# Row 41 - Aptos
Set-TextValue "D41" "13.63"
$ws.Range("E41").Value = "  +2.75%  "

# Row 42 - TrustWalletToken
Set-TextValue "D42" "1.429"
$ws.Range("E42").Value = "  -4.65%  "

# Row 43 - TheSandbox
Set-TextValue "D43" "0.7511"
$ws.Range("E43").Value = "  -3.29%  "

# Row 44 - EnergySwap
Set-TextValue "D44" "16.03"
$ws.Range("E44").Value = "  -1.22%  "

# Row 45 - Decentraland
Set-TextValue "D45" "0.6921"
$ws.Range("E45").Value = "  -1.01%  "

# Row 46 - NEARProtocol
Set-TextValue "D46" "2.471"
$ws.Range("E46").Value = "  -2.68%  "

# Row 47 - PancakeSwap
Set-TextValue "D47" "4.070"
$ws.Range("E47").Value = "  -1.88%  "

# Row 48 - Frax
$ws.Range("E48").Value = "  -0.53%  "

# Row 49 - Cronos
Set-TextValue "D49" "0.08247"
$ws.Range("E49").Value = "  -3.62%  "

# Row 50 - Quant
Set-TextValue "D50" "135.24"
$ws.Range("E50").Value = "  -1.10%  "

# Row 51 - Flow
Set-TextValue "D51" "1.200"
$ws.Range("E51").Value = "  -8.40%  "
